$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 50709.26919481565
$ws.Range("C2").Value = 32436.48727425859
$ws.Range("D2").Value = 30527.88561606283
$ws.Range("B3").Value = 32192.30262320694
$ws.Range("C3").Value = 22866.78012130396
$ws.Range("D3").Value = 24205.3779859723
$ws.Range("B5").Value = 1872.448421409234
$ws.Range("C5").Value = -5293.074080493738
$ws.Range("D5").Value = -3954.476215825402
$ws.Range("B6").Value = 3504.077841354558
$ws.Range("C6").Value = 600.0792247045911
$ws.Range("D6").Value = 851.4687391656371
$ws.Range("B8").Value = 1936.577841354558
$ws.Range("C8").Value = -967.4207752954089
$ws.Range("D8").Value = -716.0312608343629
$ws.Range("B9").Value = -13.57724609374998
$ws.Range("C9").Value = -1.139196777343727
$ws.Range("D9").Value = -0.8878845214843523
$ws.Range("B10").Value = 26.99999999997857
$ws.Range("C10").Value = 23.71977692081987
$ws.Range("D10").Value = 26.79777666228358
$ws.Range("B11").Value = 26.63446859944048
$ws.Range("C11").Value = 23.99999999999972
$ws.Range("D11").Value = 26.99999999995276
$ws.Range("B12").Value = 27.02262291110753
$ws.Range("C12").Value = 23.25632273469279
$ws.Range("D12").Value = 26.34428433441923
$ws.Range("B13").Value = 64.68334894098598
$ws.Range("C13").Value = 22.30279351510057
$ws.Range("D13").Value = 25.58199758956857
$ws.Range("B14").Value = 11.76268661558453
$ws.Range("C14").Value = 16.05967017558862
$ws.Range("D14").Value = 16.52853569592492
$ws.Range("B15").Value = 8.42226987810564
$ws.Range("C15").Value = 6.579533647300682
$ws.Range("D15").Value = 6.398559916072532
$ws.Range("B16").Value = 15.44572833419447
$ws.Range("C16").Value = 12.46304688888272
$ws.Range("D16").Value = 12.3467559851574
$ws.Range("B17").Value = 15.07037301409483
$ws.Range("C17").Value = 13.00171485776169
$ws.Range("D17").Value = 12.80614804761802
$ws.Range("B18").Value = 1.216857408920424
$ws.Range("C18").Value = 7.166476143855323
$ws.Range("D18").Value = 6.873220772578689
$ws.Range("B19").Value = 1.999999999999922
$ws.Range("C19").Value = 2.120000000000009
$ws.Range("D19").Value = 1.760000000014667
$ws.Range("B20").Value = 0.218890967250746
$ws.Range("C20").Value = 0.2446960102280806
$ws.Range("D20").Value = 0.2061137846734573
$ws.Range("B21").Value = 1.816109032749176
$ws.Range("C21").Value = 1.892303989771928
$ws.Range("D21").Value = 1.57088621534121
$ws.Range("B22").Value = 1.819299135049805
$ws.Range("C22").Value = 1.895350387741783
$ws.Range("D22").Value = 1.573887308414705
$ws.Range("B23").Value = 0.2189495116472244
$ws.Range("C23").Value = 0.2447557151317596
$ws.Range("D23").Value = 0.2061730474233627
$ws.Range("B25").Value = 1.038248658180237
$ws.Range("C25").Value = 1.080106139183044
$ws.Range("D25").Value = 0.9000603556632996
$ws.Range("B26").Value = 0.1228958829729158
$ws.Range("C26").Value = 0.1203273953121686
$ws.Range("D26").Value = 0.1203866854432739
$ws.Range("B29").Value = 104.5986171875
$ws.Range("C29").Value = 79.60111718749999
$ws.Range("D29").Value = 78.59155468749999
$ws.Range("B30").Value = 100.1430703125
$ws.Range("C30").Value = 73.89046875
$ws.Range("D30").Value = 74.28260937500001
$ws.Range("B31").Value = 104.5986147265888
$ws.Range("C31").Value = 79.60111683149745
$ws.Range("D31").Value = 78.59155668363148
